# "creato modulo per certificates"
# The workbook originally ships an empty placeholder sheet ("Sheet") plus the
# real certificates export ("Output"). This edit drops the empty placeholder
# and promotes the certificates data to be the workbook's single sheet named
# "Sheet", then backfills Issuer / Issue date / Settlement currency on the
# four rows that were missing them (21-24), matching row 20's values.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Remove the empty placeholder sheet.
$placeholder = $wb.Worksheets.Item("Sheet")
$placeholder.Delete()

# The remaining "Output" sheet becomes the workbook's single sheet, renamed
# back to "Sheet".
$ws = $wb.Worksheets.Item("Output")
$ws.Name = "Sheet"

# Backfill the four data rows (21-24) with the same Issuer / Issue date /
# Settlement currency text already present on row 20.
$issuer = "BNP Paribas Issuance B.V."
$issueDate = "26 May 2020."
$settlementCurrency = 'The settlement currency for the payment of the Cash Settlement Amount is Euro ("EUR")'

foreach ($row in 21..24) {
    $ws.Cells.Item($row, 1).Value = $issuer
    $ws.Cells.Item($row, 8).Value = $issueDate
    $ws.Cells.Item($row, 40).Value = $settlementCurrency
}
